{"js": "// Append \" (Changed main)\" as three new separate runs to the end of the\n// first paragraph (\"This is a Microsoft word document.\"), matching:\n//   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>\n//   <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n//   <w:r><w:t>Changed main</w:t></w:r>\n//   <w:r><w:t>)</w:t></w:r>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Use insertOoxml so each appended piece of text lands in its own <w:r>\n// run (matching the target diff) instead of being merged into the\n// existing run the way insertText would.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:r><w:t>Changed main</w:t></w:r>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nfirstParagraph.insertOoxml(ooxml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append \" (Changed main)\" as three new separate runs to the end of the\n# first paragraph (\"This is a Microsoft word document.\"), matching:\n#   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>\n#   <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n#   <w:r><w:t>Changed main</w:t></w:r>\n#   <w:r><w:t>)</w:t></w:r>\n\n$d = $word.ActiveDocument\n\n# Locate the end of the target sentence with Find (robust to it not being\n# the very first range) and collapse the match down to an insertion point.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"This is a Microsoft word document.\")\nif (-not $found) {\n    throw \"Could not find target paragraph text\"\n}\n$rng.Collapse(0) | Out-Null   # wdCollapseEnd\n\n# Plain Range.InsertAfter()/Range.Text would just get absorbed into the\n# neighbouring run's <w:t>, giving one merged run. To land the new text in\n# three distinct <w:r> runs (as the diff requires) give the range real\n# extent with a one-character placeholder, then use InsertXML - which\n# REPLACES a Range's contents with parsed OOXML - to splice in the exact\n# run structure we want.\n$rng.InsertAfter(\"X\") | Out-Null\n\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n<w:r><w:t>Changed main</w:t></w:r>\n<w:r><w:t>)</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$rng.InsertXML($ooxml) | Out-Null\n"}
